$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2's date (2022-06-08, stored as a date serial) becomes the plain
# number 220608 (yymmdd as digits, no longer an actual date) - drop the
# "yymmdd;@" date format on column A (A1:A3) back to a plain/general number.
$ws.Range("A1:A3").NumberFormat = "general"

$ws.Range("A2").Value = 220608

# New row 3: EffNetV2M run with Class count 11 (박영서)
$ws.Range("A3").Value = 220609
$ws.Range("B3").Value = "박영서"
$ws.Range("C3").Value = "EffNetV2M"
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 0.9192
$ws.Range("F3").Value = 0.9168
$ws.Range("G3").Value = 0.463
$ws.Range("H3").Value = 0.3789

$ws.Range("F9").Select() | Out-Null
